$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Procedure text for the "wendlandt" row (row 4) with the full/extended
# experimental write-up (previously truncated text).
$fullText = "To a solution of Wieland-Miescher ketone (1.0 g, 5.6 mmol) exhibiting (a], + 100 (c 1.0, toluene) in methanol (30 mL) at 0°C was added, dropwise, a solution of NaBH, (0.14 g, 3.7 mmol) in methanol (20 mL). The solution was stirred at 0°C for 30 min and several drops of acetic acid were added. The mixture was concentrated under vacuum, water was added, and the product was extracted into ether (3 X 20 mL), dried, concentrated, and purified by flash chromatography (elution with 30% ethyl acetate in petroleum ether) to give 0.980 g (97%) of the 2ß-ol-S-one as a colorless oil; [al + 183 (c 1.7, CHCl); ir (CHCl,): 3600, 3450, 1665, and 1620 cm'; 'H mr (300 MHz, CDCl;) &; 5.73 (s, 1H), 3.37 (dd, J = 1 1.5, 4.4 Hz, 1 H), 2.66 (br s, lH), 2.45-2.10 (series of m, 6H), 1.82 (m, 2H), 1.64 (m, IH), 1.45 (m, lH), 1.15 (s, 3H); I3c nrnr (75 MHz, CDCl,): 199.84, 169.04, 125.17, 77.95, 41.58, 34.10, 33.59, 31.96, 30.09, 23.08, 15.19 ppm; ms m/z (M+) calcd.: 180.1 150; found: 180.1 167."

$ws.Range("E4").Value = $fullText

# Touch column G width so it becomes an explicitly-sized column (matches the
# author re-saving after interacting with / auto-fitting that column).
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Move the active selection to B4, as it ended up after the edit.
$ws.Range("B4").Select()

$wb.Save()
